# Updated cryptos list (prices + 1h volume %) as captured on
# Sun Jun 16 14:28:05 UTC 2024 with GitHub Actions.
#
# D-column price cells are stored as text in the source data (they can
# contain multiple '.' separators, e.g. "66.666.24"), so for every D-cell
# we briefly force a Text number format before writing the value (so
# Excel doesn't auto-coerce plain decimals like "609.26" into a Number),
# then restore the "Normal" style afterwards so no visible formatting
# change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "66.666.24"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +0.63%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.597.14"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "609.26"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "146.56"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "3.594.25"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +0.03%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "7.95"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -2.16%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.416"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +1.12%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "4.212.02"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.0000209"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "30.06"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -0.61%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "3.582.98"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "66.775.43"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +0.75%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "15.06"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +1.03%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "433.08"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +3.12%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "79.16"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "3.739.28"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  +1.30%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "8.11"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  +1.22%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.595.87"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "25.53"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "1.45"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  -1.64%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "7.85"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  +0.03%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "1.72"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "5.64"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "174.60"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.0855"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -0.30%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "5.24"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.895"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "1.93"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +1.13%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "45.92"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.54"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("E48").Value = "  -1.53%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "25.07"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -3.79%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "23.77"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +3.50%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "7.22"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +1.09%  "
